$wb = $excel.ActiveWorkbook

# --- hunk 0: sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 108.6
$ws.Range("I4").Value = 118.5
$ws.Range("J4").Value = 69
$ws.Range("K4").Value = 118.5
$ws.Range("L4").Value = 69
$ws.Range("M4").Value = -4.5
$ws.Range("N4").Value = -297

# --- hunk 1: sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1334.2745
$ws.Range("I132").Value = 1334.6666
$ws.Range("J132").Value = 1332.4445
$ws.Range("K132").Value = 4003.9998
$ws.Range("L132").Value = 3997.3335
$ws.Range("M132").Value = -1473.9998
$ws.Range("N132").Value = -9057.333500000001

# --- hunk 2: sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1178.6774
$ws.Range("I135").Value = 1158.6666
$ws.Range("J135").Value = 1247.2858
$ws.Range("K135").Value = 10427.9994
$ws.Range("L135").Value = 11225.5722
$ws.Range("M135").Value = -7892.999400000001
$ws.Range("N135").Value = -16295.5722

# --- hunk 3: sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2844.2454
$ws.Range("I138").Value = 790.375
$ws.Range("J138").Value = 4544
$ws.Range("K138").Value = 2371.125
$ws.Range("L138").Value = 13632
$ws.Range("M138").Value = 2768.875
$ws.Range("N138").Value = -23912

# --- hunk 4: sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2759.9048
$ws.Range("I141").Value = 2637.2632
$ws.Range("J141").Value = 3925
$ws.Range("K141").Value = 7911.7896
$ws.Range("L141").Value = 11775
$ws.Range("M141").Value = -2731.7896
$ws.Range("N141").Value = -22135

# --- hunk 5: sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 300
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -532

# --- hunk 6: sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4406.8887
$ws.Range("I61").Value = 2665.3403
$ws.Range("K61").Value = 2665.3403
$ws.Range("M61").Value = -2453.3403

# --- hunk 7: sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7467
$ws.Range("I74").Value = 8461.096
$ws.Range("J74").Value = 3671.3635
$ws.Range("K74").Value = 8461.096
$ws.Range("L74").Value = 3671.3635
$ws.Range("M74").Value = -7587.096
$ws.Range("N74").Value = -5419.363499999999

# --- hunk 8: sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

# --- hunk 9: sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 7467
$ws.Range("I77").Value = 8461.096
$ws.Range("J77").Value = 3671.3635
$ws.Range("K77").Value = 42305.48
$ws.Range("L77").Value = 18356.8175
$ws.Range("M77").Value = -37937.48
$ws.Range("N77").Value = -27092.8175

# --- hunk 10: sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

# --- hunk 11: sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 10400
$ws.Range("J124").Value = 10400
$ws.Range("L124").Value = 10400
$ws.Range("N124").Value = -20220

# --- hunk 12: sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H131").Value = 40715
$ws.Range("J131").Value = 40715
$ws.Range("L131").Value = 40715
$ws.Range("N131").Value = -50795

# --- hunk 13: sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4406.8887
$ws.Range("I136").Value = 2665.3403
$ws.Range("K136").Value = 7996.0209
$ws.Range("M136").Value = -5446.0209

# --- hunk 14: sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# --- hunk 15: sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5209.4126
$ws.Range("I134").Value = 2527.8518
$ws.Range("J134").Value = 21298.777
$ws.Range("K134").Value = 7583.555399999999
$ws.Range("L134").Value = 63896.33099999999
$ws.Range("M134").Value = -5048.555399999999
$ws.Range("N134").Value = -68966.33099999999

# --- hunk 16: sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 484.81482
$ws.Range("I7").Value = 449.68182
$ws.Range("K7").Value = 449.68182
$ws.Range("M7").Value = -336.68182

# --- hunk 17: sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4425.778
$ws.Range("I31").Value = 3871.7334
$ws.Range("J31").Value = 5118.3335
$ws.Range("K31").Value = 3871.7334
$ws.Range("L31").Value = 5118.3335
$ws.Range("M31").Value = -3576.7334
$ws.Range("N31").Value = -5708.3335

# --- hunk 18: sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4425.778
$ws.Range("I34").Value = 3871.7334
$ws.Range("J34").Value = 5118.3335
$ws.Range("K34").Value = 3871.7334
$ws.Range("L34").Value = 5118.3335
$ws.Range("M34").Value = -3669.7334
$ws.Range("N34").Value = -5522.3335

# --- hunk 19: sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1703
$ws.Range("I105").Value = 1999
$ws.Range("K105").Value = 1999
$ws.Range("M105").Value = -252

# --- hunk 20: sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 19092.193
$ws.Range("I132").Value = 11280.55
$ws.Range("J132").Value = 39470.39
$ws.Range("K132").Value = 33841.64999999999
$ws.Range("L132").Value = 118411.17
$ws.Range("M132").Value = -31311.64999999999
$ws.Range("N132").Value = -123471.17

# --- hunk 21: sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2427.14
$ws.Range("I134").Value = 1964.5349
$ws.Range("J134").Value = 5268.857
$ws.Range("K134").Value = 5893.6047
$ws.Range("L134").Value = 15806.571
$ws.Range("M134").Value = -3358.6047
$ws.Range("N134").Value = -20876.571

# --- hunk 22: sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 5400
$ws.Range("I137").Value = 5000
$ws.Range("J137").Value = 6000
$ws.Range("K137").Value = 15000
$ws.Range("L137").Value = 18000
$ws.Range("M137").Value = -9900
$ws.Range("N137").Value = -28200

# --- hunk 23: sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 5457.8823
$ws.Range("I140").Value = 2889.5
$ws.Range("J140").Value = 9127
$ws.Range("K140").Value = 8668.5
$ws.Range("L140").Value = 27381
$ws.Range("M140").Value = -3488.5
$ws.Range("N140").Value = -37741

# --- hunk 24: sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2419.2727
$ws.Range("I102").Value = 2541.2
$ws.Range("J102").Value = 1200
$ws.Range("K102").Value = 2541.2
$ws.Range("L102").Value = 1200
$ws.Range("M102").Value = -919.1999999999998
$ws.Range("N102").Value = -4444

# --- hunk 25: sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9282
$ws.Range("I132").Value = 6919.436
$ws.Range("J132").Value = 20799.5
$ws.Range("K132").Value = 20758.308
$ws.Range("L132").Value = 62398.5
$ws.Range("M132").Value = -18228.308
$ws.Range("N132").Value = -67458.5

# --- hunk 26: sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5889.174
$ws.Range("I132").Value = 4211.4443
$ws.Range("J132").Value = 6967.7144
$ws.Range("K132").Value = 12634.3329
$ws.Range("L132").Value = 20903.1432
$ws.Range("M132").Value = -10104.3329
$ws.Range("N132").Value = -25963.1432

# --- hunk 27: sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4035.0488
$ws.Range("I136").Value = 3339.303
$ws.Range("J136").Value = 6905
$ws.Range("K136").Value = 10017.909
$ws.Range("L136").Value = 20715
$ws.Range("M136").Value = -7467.909
$ws.Range("N136").Value = -25815

# --- hunk 28: sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 14555.052
$ws.Range("I132").Value = 8508.768
$ws.Range("J132").Value = 29945.592
$ws.Range("K132").Value = 25526.304
$ws.Range("L132").Value = 89836.776
$ws.Range("M132").Value = -22996.304
$ws.Range("N132").Value = -94896.776

# --- hunk 29: sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4001864
$ws.Range("I136").Value = 5264167
$ws.Range("K136").Value = 15792501
$ws.Range("M136").Value = -15789951
